# Replicape BOM: add "Alt:" beaglebone-LCD3-cape line, fix "NONE" -> "None"
# (see commit "Committing before going to 3DP weekend : )")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "Optional:" kit row (row 38) had a typo'd "NONE" in its
# alternate-part column; the shared string backing C38 is reused by the new
# row below, so fix it in place so both cells end up reading "None".
$ws.Range("C38").Value = "None"

# Row 40 used to hold a lone "Alt:" label. Drop it - it gets replaced by a
# fully populated row one line up (row 39).
[void]$ws.Rows.Item(40).Delete()

# New row 39: the alternate BeagleBone LCD3 cape option.
$ws.Range("A39").Value = "BB-BONE-LCD3-01-ND"
$ws.Range("B39").Value = "BEAGLEBONE LCD3 CAPE"
$ws.Range("C39").Value = "None"
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 69.95

# Part nr. / description cells wrap instead of overflowing.
$ws.Range("A39:B39").WrapText = $true

# Row needs to be a touch taller to fit the wrapped text.
$ws.Rows.Item(39).RowHeight = 14.9

# Column A widened slightly to accommodate the new, longer part numbers.
$ws.Columns.Item(1).ColumnWidth = 18.3

# Restore the cursor to where it was left (one row below the new content).
[void]$ws.Range("C40").Select()
